$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot full rows 2-45 (columns A:T) before rewriting, since this edit
# is a reordering (permutation) of existing records by date.
$snapshot = @{}
for ($r = 2; $r -le 45; $r++) {
    $rowVals = $ws.Range("A$r`:T$r").Value2
    $snapshot[$r] = $rowVals
}

# Write rows back out in their new order (after-row <- before-row).
$ws.Range("A2`:T2").Value2 = $snapshot[32]
$ws.Range("A3`:T3").Value2 = $snapshot[33]
$ws.Range("A4`:T4").Value2 = $snapshot[11]
$ws.Range("A5`:T5").Value2 = $snapshot[12]
$ws.Range("A6`:T6").Value2 = $snapshot[40]
$ws.Range("A7`:T7").Value2 = $snapshot[41]
$ws.Range("A8`:T8").Value2 = $snapshot[45]
$ws.Range("A9`:T9").Value2 = $snapshot[2]
$ws.Range("A10`:T10").Value2 = $snapshot[5]
$ws.Range("A11`:T11").Value2 = $snapshot[13]
$ws.Range("A12`:T12").Value2 = $snapshot[14]
$ws.Range("A13`:T13").Value2 = $snapshot[8]
$ws.Range("A14`:T14").Value2 = $snapshot[9]
$ws.Range("A15`:T15").Value2 = $snapshot[10]
$ws.Range("A16`:T16").Value2 = $snapshot[18]
$ws.Range("A17`:T17").Value2 = $snapshot[19]
$ws.Range("A18`:T18").Value2 = $snapshot[22]
$ws.Range("A19`:T19").Value2 = $snapshot[23]
$ws.Range("A20`:T20").Value2 = $snapshot[44]
$ws.Range("A21`:T21").Value2 = $snapshot[39]
$ws.Range("A22`:T22").Value2 = $snapshot[36]
$ws.Range("A23`:T23").Value2 = $snapshot[20]
$ws.Range("A24`:T24").Value2 = $snapshot[21]
$ws.Range("A25`:T25").Value2 = $snapshot[26]
$ws.Range("A26`:T26").Value2 = $snapshot[27]
$ws.Range("A27`:T27").Value2 = $snapshot[31]
$ws.Range("A28`:T28").Value2 = $snapshot[6]
$ws.Range("A29`:T29").Value2 = $snapshot[7]
$ws.Range("A30`:T30").Value2 = $snapshot[35]
$ws.Range("A31`:T31").Value2 = $snapshot[15]
$ws.Range("A32`:T32").Value2 = $snapshot[29]
$ws.Range("A33`:T33").Value2 = $snapshot[30]
$ws.Range("A34`:T34").Value2 = $snapshot[24]
$ws.Range("A35`:T35").Value2 = $snapshot[25]
$ws.Range("A36`:T36").Value2 = $snapshot[28]
$ws.Range("A37`:T37").Value2 = $snapshot[34]
$ws.Range("A38`:T38").Value2 = $snapshot[42]
$ws.Range("A39`:T39").Value2 = $snapshot[43]
$ws.Range("A40`:T40").Value2 = $snapshot[16]
$ws.Range("A41`:T41").Value2 = $snapshot[17]
$ws.Range("A42`:T42").Value2 = $snapshot[37]
$ws.Range("A43`:T43").Value2 = $snapshot[38]
$ws.Range("A44`:T44").Value2 = $snapshot[3]
$ws.Range("A45`:T45").Value2 = $snapshot[4]

Write-Host "done"
